# Updates the cryptocurrency price/volume table on Sheet1 (rows 2-51)
# to reflect the latest scraped values, as published by the GitHub
# Actions job that refreshes cryptos.xlsx.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> updated "Price" (column D) and/or "Volume(1h)" (column E) values.
# Rows without a "D" entry only had their Volume(1h) percentage updated.
$updates = @{
    2 = @{ D = "60.830.78"; E = "  -1.51%  " }
    3 = @{ D = "2.910.18"; E = "  -2.38%  " }
    4 = @{ E = "  +0.02%  " }
    5 = @{ D = "527.89"; E = "  -2.36%  " }
    6 = @{ D = "144.04"; E = "  -5.13%  " }
    7 = @{ D = "0.998"; E = "  -0.04%  " }
    8 = @{ D = "0.547"; E = "  -3.67%  " }
    9 = @{ D = "2.922.34"; E = "  -2.34%  " }
    10 = @{ E = "  -4.32%  " }
    11 = @{ D = "6.06"; E = "  -1.53%  " }
    12 = @{ D = "0.359"; E = "  -2.64%  " }
    13 = @{ D = "3.415.82"; E = "  -2.44%  " }
    14 = @{ E = "  +3.23%  " }
    15 = @{ D = "60.776.89"; E = "  -1.66%  " }
    16 = @{ D = "22.58"; E = "  -5.60%  " }
    17 = @{ D = "2.922.18"; E = "  -2.15%  " }
    18 = @{ D = "0.0000141"; E = "  -4.03%  " }
    19 = @{ D = "4.97"; E = "  -3.75%  " }
    20 = @{ D = "11.61"; E = "  -3.48%  " }
    21 = @{ D = "354.93"; E = "  -6.85%  " }
    22 = @{ E = "  -2.66%  " }
    23 = @{ E = "  +0.08%  " }
    24 = @{ E = "  +1.30%  " }
    25 = @{ D = "65.07"; E = "  -1.39%  " }
    26 = @{ D = "0.452"; E = "  -3.96%  " }
    27 = @{ E = "  -6.30%  " }
    28 = @{ E = "  -0.10%  " }
    29 = @{ E = "  -2.97%  " }
    30 = @{ D = "0.0₃0858"; E = "  -8.77%  " }
    32 = @{ D = "1.69"; E = "  -1.50%  " }
    33 = @{ D = "19.64"; E = "  -4.06%  " }
    34 = @{ D = "153.90"; E = "  -3.69%  " }
    35 = @{ D = "4.39"; E = "  -3.98%  " }
    36 = @{ D = "5.58"; E = "  -5.81%  " }
    37 = @{ D = "0.998"; E = "  -6.88%  " }
    38 = @{ D = "1.20"; E = "  -5.56%  " }
    39 = @{ D = "37.53"; E = "  -0.13%  " }
    40 = @{ E = "  -4.57%  " }
    41 = @{ E = "  -4.40%  " }
    42 = @{ D = "2.293.69"; E = "  -5.30%  " }
    43 = @{ D = "0.652"; E = "  -2.90%  " }
    44 = @{ D = "0.0583"; E = "  -0.95%  " }
    45 = @{ D = "20.41"; E = "  -7.17%  " }
    46 = @{ E = "  +0.08%  " }
    47 = @{ D = "4.98"; E = "  -3.72%  " }
    48 = @{ E = "  -2.79%  " }
    49 = @{ E = "  -0.72%  " }
    50 = @{ D = "0.0918"; E = "  -3.56%  " }
    51 = @{ D = "18.45"; E = "  -6.88%  " }
}

foreach ($rowNum in $updates.Keys) {
    $item = $updates[$rowNum]

    if ($item.ContainsKey("D")) {
        $cell = $ws.Cells.Item($rowNum, 4)
        # Force text format so values like "60.830.78" or "1.20" are kept
        # as literal strings instead of being re-interpreted as numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $item.D
    }

    if ($item.ContainsKey("E")) {
        $ws.Cells.Item($rowNum, 5).Value = $item.E
    }
}
